$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "Nos" (nodes: x, y coordinates) -> grows from 4 to 9 data rows
# ---------------------------------------------------------------
$wsNos = $wb.Worksheets.Item("Nos")
$wsNos.Range("A3").Value = 4
$wsNos.Range("B3").Value = 4
$wsNos.Range("A4").Value = 8
$wsNos.Range("B4").Value = 8

# Extend the existing "A5:B5" formatting (centred, style used for data rows)
# down through row 10 before filling in the new coordinate pairs.
$wsNos.Range("A5:B5").Copy()
$wsNos.Range("A5:B10").PasteSpecial(-4122)
$wsNos.Range("A5").Value = 4
$wsNos.Range("B5").Value = 0
$wsNos.Range("A6").Value = 8
$wsNos.Range("B6").Value = 4
$wsNos.Range("A7").Value = 8
$wsNos.Range("B7").Value = 0
$wsNos.Range("A8").Value = 12
$wsNos.Range("B8").Value = 4
$wsNos.Range("A9").Value = 12
$wsNos.Range("B9").Value = 0
$wsNos.Range("A10").Value = 16
$wsNos.Range("B10").Value = 0
$wsNos.Range("B10").Select()

# ---------------------------------------------------------------
# Sheet "Incidencia" (member connectivity: node1, node2, E, A) -> grows
# from 3 to 16 data rows
# ---------------------------------------------------------------
$wsInc = $wb.Worksheets.Item("Incidencia")
$wsInc.Range("A3").Value = 1
$wsInc.Range("B3").Value = 4
$wsInc.Range("A4").Value = 2
$wsInc.Range("B4").Value = 4

# Extend row 4's look (E, A columns formatted) down through row 17, then
# populate the new connectivity pairs with matching E / A values.
$wsInc.Range("A4:D4").Copy()
$wsInc.Range("A5:D17").PasteSpecial(-4122)

$incData = @(
  @(2, 3),
  @(2, 5),
  @(3, 5),
  @(3, 7),
  @(4, 6),
  @(4, 5),
  @(5, 6),
  @(5, 7),
  @(5, 8),
  @(6, 8),
  @(7, 8),
  @(7, 9),
  @(8, 9)
)
$r = 5
foreach ($pair in $incData) {
  $wsInc.Range("A$r").Value = $pair[0]
  $wsInc.Range("B$r").Value = $pair[1]
  $wsInc.Range("C$r").Value = 210000000000
  $wsInc.Range("D$r").Value = 0.0002
  $r = $r + 1
}
$wsInc.Range("D17").Select()

# ---------------------------------------------------------------
# Sheet "Carregamento" (loads: node, direction, load [N]) -> one more row
# ---------------------------------------------------------------
$wsCar = $wb.Worksheets.Item("Carregamento")
$wsCar.Range("C2").Value = 200
$wsCar.Range("C3").Value = -1000
$wsCar.Range("A4").Value = 2
$wsCar.Range("B4").Value = 1
$wsCar.Range("C4").Value = 200
$wsCar.Range("C5").Select()

# ---------------------------------------------------------------
# Sheet "Restricao" (restraints: node, direction) -> grows from 3 to 7 rows
# ---------------------------------------------------------------
$wsRes = $wb.Worksheets.Item("Restricao")
$wsRes.Range("A3").Value = 1
$wsRes.Range("B3").Value = 2
$wsRes.Range("A4").Value = 4

# Row 5 inherits the centred style already used by row 5's blank cells
# (same as "Nos"); extend rows 6-8 with the plain style used by rows 6-8
# previously (same as row 4 here).
$wsRes.Range("A5").Value = 6
$wsRes.Range("B5").Value = 2

$wsRes.Range("A4:B4").Copy()
$wsRes.Range("A6:B8").PasteSpecial(-4122)
$wsRes.Range("A6").Value = 8
$wsRes.Range("B6").Value = 2
$wsRes.Range("A7").Value = 9
$wsRes.Range("B7").Value = 1
$wsRes.Range("A8").Value = 9
$wsRes.Range("B8").Value = 2

$wsRes.Range("H37").Select()
